# Research Question edit: remove "departing from Delhi" / fix spacing
# so the question covers the whole dataset instead of just Delhi flights.

$p = $ppt.ActivePresentation

function Replace-InTextRange {
    param($TextRange, $Find, $Replacement)
    $full = $TextRange.Text
    $idx = $full.IndexOf($Find)
    if ($idx -lt 0) {
        return $false
    }
    $sub = $TextRange.Characters($idx + 1, $Find.Length)
    $sub.Text = $Replacement
    return $true
}

# --- Slide 3: "Is there a correlation ... domestic flights departing from Delhi in India? ." ---
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item("Title 4").TextFrame.TextRange

Replace-InTextRange $title3 "domestic flights departing from Delhi " "domestic flights " | Out-Null
Replace-InTextRange $title3 "in India? ." "in India?." | Out-Null

# --- Slide 4: Null / Alternative hypothesis wording ---
$slide4 = $p.Slides.Item(4)
$title4 = $slide4.Shapes.Item("Title 8").TextFrame.TextRange

Replace-InTextRange $title4 " for domestic flights departing from Delhi in India" " for domestic flights departing in India" | Out-Null
Replace-InTextRange $title4 "for domestic flights departing from Delhi in India." "for domestic flights in India." | Out-Null
